# Simplify the Workblocks sheet: drop the per-workblock "SuppressSuccessful"
# rows (exceptions are now always rethrown, so there is nothing to suppress),
# rename the Init workblock to "InitAllApplications" and add a new
# "CloseAllApplications" workblock.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Workblocks")

# Row 2: wbInit_Type -> wbInitAllApplications_Type (value/description unchanged)
$ws.Cells.Item(2, 1).Value = "wbInitAllApplications_Type"

# Row 3 used to be wbInit_SuppressSuccessful; reuse it for
# wbGetTransactionData_Type (previously row 4).
$ws.Cells.Item(3, 1).Value = "wbGetTransactionData_Type"
$ws.Cells.Item(3, 2).Value = "GetData"
$ws.Cells.Item(3, 3).Value = "Name of Workblock"

# Row 4 used to be wbGetTransactionData_SuppressSuccessful; reuse it for
# wbProcessTransaction_Type (previously row 6).
$ws.Cells.Item(4, 1).Value = "wbProcessTransaction_Type"
$ws.Cells.Item(4, 2).Value = "Process"
$ws.Cells.Item(4, 3).Value = "Name of Workblock"

# Row 5 used to be wbProcessTransaction_Type; reuse it for the new
# wbCloseAllApplications_Type workblock.
$ws.Cells.Item(5, 1).Value = "wbCloseAllApplications_Type"
$ws.Cells.Item(5, 2).Value = "Close"
$ws.Cells.Item(5, 3).Value = "Name of Workblock"

# The old rows 6 (wbProcessTransaction_Type) and 7
# (wbProcessTransaction_SuppressSuccessful) have already been folded into
# rows 2-5 above, so drop the now-redundant trailing rows (delete bottom-up
# so row numbers stay stable).
$ws.Rows.Item(7).Delete()
$ws.Rows.Item(6).Delete()

# The workbook now opens on the Workblocks tab instead of Credentials.
$ws.Activate() | Out-Null
$ws.Range("B11").Select() | Out-Null
